$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "first row" manual-seed values (A4:C4): new current price, new
# dust-amount and a refreshed date.
$ws.Range("A4").Value = 999
$ws.Range("B4").Value = 0.0000000001
$ws.Range("C4").Value = 45972

# Extend the same formatting that A4:C4 already carries down into the next
# couple of rows (A5:C6) so newly-entered purchase rows line up visually.
$ws.Range("A4:B4").Copy()
$ws.Range("A5:B6").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4122)

# Move the active selection down to where the next rows of data will go.
[void]$ws.Range("A5:D19").Select()
$excel.CutCopyMode = 0

[void]$wb.Application.Calculate()
